# PCA + rbf SVM
#
# 1) Inserts a new results row (SVM (SVC) (leader) / Rbf, C=1, gamma=0.0001)
#    right after the existing "Linear, C=1" row, pushing the rest of the
#    comparison table down by one row.
# 2) Appends a brand-new row at the bottom of the table for the
#    PCA(26, true) + rbf-SVM result.
# 3) Widens column B so the longer parameter strings fit, and leaves the
#    selection where the author left off (B12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row above the current "kNN" row (row 3) -----------------
# This shifts the existing rows 3..6 down to 4..7.
$ws.Rows.Item(3).Insert()

# The freshly inserted row loses the bordered "data row" style, so pull the
# formatting back from row 2 (same style as every other data row) before
# filling in values.
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A3:G3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "SVM (SVC) (лидер)"
$ws.Range("B3").Value = "Rbf, C=1, gamma=0.0001"
$ws.Range("C3").Value = 5000
$ws.Range("D3").Value = "CV, 5"
$ws.Range("E3").Value = 0.78
$ws.Range("F3").Value = 0.02
$ws.Range("G3").Value = 10.199999999999999

# --- Append a brand-new row for the PCA+rbf SVM result ----------------------
# After the insert-shift above, the table's last populated row is row 6
# ("RandomizedPCA+SVM" / "PCA(27, true), SVM(Linear, C=1)"); copy its format
# down onto the new row 7.
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A7:G7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A7").Value = "RandomizedPCA+SVM"
$ws.Range("B7").Value = "PCA(26, true), SVM(rbf, C=1, gamma=0,0001)"
$ws.Range("C7").Value = 5000
$ws.Range("D7").Value = "CV, 5"
$ws.Range("E7").Value = 0.66
$ws.Range("F7").Value = 0.01
$ws.Range("G7").Value = 4.8099999999999996

# --- Cosmetic follow-ups from the authored diff -----------------------------
# Column B got wider to fit the longer parameter strings (~40.71 characters;
# Excel quantizes column widths to whole pixels, so 39.8 is the input that
# lands on the closest achievable stored width).
$ws.Columns.Item(2).ColumnWidth = 39.8

# Selection moved to B12 (author was about to type something below the table).
$ws.Range("B12").Select()
